# Trade #9 closed at 2026-02-17 19:44:47 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1300      # Current Capital
$summary.Range("B4").Value = 0         # Total P&L $
$summary.Range("B5").Value = 0         # Total P&L %
$summary.Range("B6").Value = 9         # Total Trades
$summary.Range("B7").Value = 5         # Winning Trades
$summary.Range("B9").Value = 55.56     # Win Rate %

# --- Sheet: Strategy Status (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100        # Capital
$status.Range("D4").Value = 9          # Trades
$status.Range("E4").Value = 0          # P&L $
$status.Range("F4").Value = 0          # P&L %
$status.Range("G4").Value = 55.56      # Win Rate %

# --- Append new closed trade (#9) to "All Trades" and "MarketMaking" sheets ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Force the date/time-looking strings to stay text (avoid Excel's
    # autoconvert-to-date behavior), then strip the format again so no
    # stray number format / style gets attached to the cell.
    $ws.Range("B10:C10").NumberFormat = "@"

    $ws.Cells.Item(10, 1).Value  = 9
    $ws.Cells.Item(10, 2).Value  = "2026-02-17"
    $ws.Cells.Item(10, 3).Value  = "19:44:40"
    $ws.Cells.Item(10, 4).Value  = "MarketMaking"
    $ws.Cells.Item(10, 5).Value  = "DOWN"
    $ws.Cells.Item(10, 6).Value  = 0.97
    $ws.Cells.Item(10, 7).Value  = 0.98
    $ws.Cells.Item(10, 8).Value  = "CLOSED"
    $ws.Cells.Item(10, 9).Value  = 1.0309
    $ws.Cells.Item(10, 10).Value = 0.01
    $ws.Cells.Item(10, 11).Value = 100
    $ws.Cells.Item(10, 12).Value = 0
    $ws.Cells.Item(10, 13).Value = 0
    $ws.Cells.Item(10, 14).Value = 0.6
    $ws.Cells.Item(10, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(10, 16).Value = "early_exit"
    $ws.Cells.Item(10, 17).Value = 0.13

    $ws.Range("B10:C10").ClearFormats()
}
